# Updated cryptos list on Thu Jul 27 18:41:54 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.175.49"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.53"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7133"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.26"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07755"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3066"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.86"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08240"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.00"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7145"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.199"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.03"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.233.26"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.822"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.03"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007763"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.128.33"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.921"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1581"
$ws.Range("E25").Value = "  +8.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.81"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.878"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.491"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.331"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05168"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.910"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.173"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7266"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.690"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.160.62"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9006"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.087"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.92"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.67"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.020.37"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5288"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.758"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.221"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.858"
$ws.Range("E50").Value = "  +3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.60%  "
